# edit.ps1
# Applies the diff: merges two runs around "Harry", inserts a
# <w:lastRenderedPageBreak/> marker, merges the "C"/"reate..." runs and
# drops the _GoBack bookmark there, then appends the shoe-tying
# instructions paragraphs (re-adding the _GoBack bookmark in its new spot).

$d = $word.ActiveDocument
$q = [char]34

# ------------------------------------------------------------------
# Change 1: 'if myname == "Harry"' -- merge the two trailing runs into one.
# ------------------------------------------------------------------
$rng1 = $d.Content
$searchText1 = "if myname == " + $q + "Harry" + $q
$found1 = $rng1.Find.Execute($searchText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $para1 = $rng1.Paragraphs(1).Range
    $para1NoMark = $d.Range($para1.Start, $para1.End - 1)
    $xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="360"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>myname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> == &quot;Harry&quot;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para1NoMark.InsertXML($xml1)
}

# ------------------------------------------------------------------
# Change 2: add <w:lastRenderedPageBreak/> before "Create a list of the
#           names ... conditional loop ..." run.
# ------------------------------------------------------------------
$rng2 = $d.Content
$searchText2 = "Create a list of the names of at least 5 of your friends and use a conditional loop to print out their names as follows:"
$found2 = $rng2.Find.Execute($searchText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $rng2.Paragraphs(1).Range
    $para2NoMark = $d.Range($para2.Start, $para2.End - 1)
    $xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:ind w:left="360"/><w:rPr><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Create a list of the names of at least 5 of your friends and use a conditional loop to print out their names as follows:</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para2NoMark.InsertXML($xml2)
}

# ------------------------------------------------------------------
# Change 3: 'Create a list of instructions for tying your shoes.' -- merge
#           the 'C' / 'reate...' runs and drop the _GoBack bookmark here
#           (it is re-added further down, in the new numbered step 10).
# ------------------------------------------------------------------
$rng3 = $d.Content
$searchText3 = "Create a list of instructions for tying your shoes."
$found3 = $rng3.Find.Execute($searchText3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $para3 = $rng3.Paragraphs(1).Range
    $para3NoMark = $d.Range($para3.Start, $para3.End - 1)
    $xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr><w:ind w:left="720"/><w:rPr><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/></w:rPr><w:t>Create a list of instructions for tying your shoes.</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para3NoMark.InsertXML($xml3)
}

# ------------------------------------------------------------------
# Change 4: insert the full shoe-tying step-by-step instructions as new
#           paragraphs right before the trailing empty ListParagraph
#           paragraph (immediately after "Provide your shoe tying
#           function below.").
# ------------------------------------------------------------------
$rng4 = $d.Content
$searchText4 = "Provide your shoe tying function below."
$found4 = $rng4.Find.Execute($searchText4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $para4 = $rng4.Paragraphs(1).Range
    $insertPoint = $d.Range($para4.End, $para4.End)
    $xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">1. </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Take 1 shoelace (the blue one) in the left hand</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> and the right shoelace (the blue</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> one)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">n the right hand. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>2. Place the red lace over the blue lace, and drop the red</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> lace.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>3. Then pick up the blue</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> lace </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">with your left hand and the red lace with your right hand (the </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>laces should be placed i</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>n a ‘X’ shape).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">4. Then put the blue lace in the opening under </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>the intersection of the two laces</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">, then drop it. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">5. Now pick the blue lace back up and pull </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>both laces tight in o</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">pposite directions. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>6. Fold the blue</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> lac</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">e in your left hand in half and </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">wrap the blue lace around you thumb. </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">7. </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Now push the b</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">lue lace through the hole after </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">you pull your thumb out with your right index finger. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">8. </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>After that pull the tip of the</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> folds </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>in the laces.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">9. </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">After that pull the bows tight. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">10. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>Now you should have a nice knot on your</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:b/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>shoe.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="720"/><w:rPr><w:sz w:val="22"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xml4)
}
